$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need a text number-format
# round-trip so Excel COM does not silently coerce the literal string into a
# floating point value (e.g. "314.30" -> 314.3000000000001). We briefly force
# a text format, assign the value, then clear the format again so the cell
# keeps its original (unstyled) appearance.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range('D2').Value = '41.570.03'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '2.466.71'
$ws.Range('E3').Value = '  -0.17%  '
Set-TextValue $ws.Range('D4') '0.998'
$ws.Range('E4').Value = '  -0.77%  '
Set-TextValue $ws.Range('D5') '314.30'
$ws.Range('E5').Value = '  +0.67%  '
Set-TextValue $ws.Range('D6') '91.21'
$ws.Range('E6').Value = '  -0.18%  '
Set-TextValue $ws.Range('D7') '0.548'
$ws.Range('E7').Value = '  +1.36%  '
$ws.Range('E8').Value = '  -0.66%  '
Set-TextValue $ws.Range('D9') '0.511'
$ws.Range('E9').Value = '  +4.37%  '
$ws.Range('E10').Value = '  -0.58%  '
Set-TextValue $ws.Range('D11') '0.0794'
$ws.Range('E11').Value = '  +2.48%  '
$ws.Range('E12').Value = '  +0.68%  '
$ws.Range('D13').Value = '2.847.12'
Set-TextValue $ws.Range('D14') '6.88'
$ws.Range('E14').Value = '  +0.91%  '
Set-TextValue $ws.Range('D15') '15.83'
$ws.Range('E15').Value = '  +4.25%  '
$ws.Range('D16').Value = '2.503.43'
$ws.Range('E16').Value = '  -1.49%  '
Set-TextValue $ws.Range('D17') '0.776'
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('D18').Value = '41.570.91'
$ws.Range('E18').Value = '  +0.63%  '
$ws.Range('E19').Value = '  +3.88%  '
$ws.Range('E20').Value = '  +2.95%  '
Set-TextValue $ws.Range('D21') '71.13'
$ws.Range('E21').Value = '  +0.88%  '
Set-TextValue $ws.Range('D22') '11.25'
$ws.Range('E22').Value = '  +3.09%  '
Set-TextValue $ws.Range('D23') '238.28'
$ws.Range('E23').Value = '  +1.65%  '
$ws.Range('E24').Value = '  +0.38%  '
$ws.Range('E25').Value = '  +2.51%  '
$ws.Range('E26').Value = '  -0.23%  '
Set-TextValue $ws.Range('D27') '24.61'
$ws.Range('E27').Value = '  +2.90%  '
$ws.Range('E28').Value = '  +0.37%  '
Set-TextValue $ws.Range('D29') '9.68'
$ws.Range('E29').Value = '  +0.28%  '
Set-TextValue $ws.Range('D30') '35.32'
$ws.Range('E30').Value = '  -1.35%  '
Set-TextValue $ws.Range('D31') '156.11'
$ws.Range('E31').Value = '  +2.60%  '
Set-TextValue $ws.Range('D32') '5.44'
$ws.Range('E32').Value = '  +0.76%  '
$ws.Range('E33').Value = '  +0.76%  '
Set-TextValue $ws.Range('D34') '0.0757'
$ws.Range('E34').Value = '  +0.82%  '
Set-TextValue $ws.Range('D35') '17.18'
$ws.Range('E35').Value = '  -0.51%  '
$ws.Range('E36').Value = '  -8.89%  '
Set-TextValue $ws.Range('D37') '2.87'
$ws.Range('E37').Value = '  -3.26%  '
$ws.Range('E38').Value = '  +1.66%  '
Set-TextValue $ws.Range('D39') '0.102'
$ws.Range('E39').Value = '  +3.55%  '
Set-TextValue $ws.Range('D40') '1.78'
$ws.Range('E40').Value = '  -2.50%  '
$ws.Range('E41').Value = '  -0.59%  '
$ws.Range('E42').Value = '  -0.88%  '
$ws.Range('D43').Value = '1.961.61'
$ws.Range('E43').Value = '  +0.03%  '
Set-TextValue $ws.Range('D45') '18.65'
$ws.Range('E45').Value = '  -2.68%  '
$ws.Range('E46').Value = '  -0.41%  '
Set-TextValue $ws.Range('D47') '9.03'
$ws.Range('E47').Value = '  +4.95%  '
$ws.Range('D48').Value = '2.706.15'
$ws.Range('E48').Value = '  -0.16%  '
Set-TextValue $ws.Range('D49') '96.90'
$ws.Range('E49').Value = '  +1.61%  '
Set-TextValue $ws.Range('D50') '67.25'
$ws.Range('E50').Value = '  -0.65%  '
$ws.Range('E51').Value = '  -1.39%  '
